# Update crypto price/volume data per latest scrape (GitHub Actions bot)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "56.729.60"
$ws.Cells.Item(3, 4).Value = "2.339.14"
$ws.Cells.Item(3, 5).Value = "  -0.25%  "
$ws.Cells.Item(4, 5).Value = "  +0.02%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "514.84"
$ws.Cells.Item(5, 5).Value = "  -0.20%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "133.78"
$ws.Cells.Item(6, 5).Value = "  +0.28%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.999"
$ws.Cells.Item(7, 5).Value = "  -0.16%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.533"
$ws.Cells.Item(8, 5).Value = "  -0.10%  "
$ws.Cells.Item(9, 5).Value = "  -1.38%  "
$ws.Cells.Item(10, 5).Value = "  -0.89%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "5.31"
$ws.Cells.Item(11, 5).Value = "  +1.39%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.339"
$ws.Cells.Item(12, 5).Value = "  +0.10%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "23.85"
$ws.Cells.Item(13, 5).Value = "  +0.98%  "
$ws.Cells.Item(14, 4).Value = "2.754.02"
$ws.Cells.Item(14, 5).Value = "  -0.16%  "
$ws.Cells.Item(15, 4).Value = "56.678.87"
$ws.Cells.Item(16, 5).Value = "  -0.13%  "
$ws.Cells.Item(17, 4).Value = "2.334.64"
$ws.Cells.Item(17, 5).Value = "  -0.46%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "10.42"
$ws.Cells.Item(18, 5).Value = "  +0.54%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "325.99"
$ws.Cells.Item(19, 5).Value = "  +2.33%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "4.18"
$ws.Cells.Item(20, 5).Value = "  -1.19%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.67"
$ws.Cells.Item(21, 5).Value = "  +1.07%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.999"
$ws.Cells.Item(22, 5).Value = "  -0.02%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "61.18"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "8.67"
$ws.Cells.Item(24, 5).Value = "  +12.38%  "
$ws.Cells.Item(25, 5).Value = "  +3.82%  "
$ws.Cells.Item(26, 5).Value = "  -0.18%  "
$ws.Cells.Item(27, 5).Value = "  +7.04%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "168.41"
$ws.Cells.Item(28, 5).Value = "  -1.32%  "
$ws.Cells.Item(29, 4).Value = "0.0₃0727"
$ws.Cells.Item(29, 5).Value = "  -1.04%  "
$ws.Cells.Item(30, 5).Value = "  +0.65%  "
$ws.Cells.Item(31, 5).Value = "  -0.92%  "
$ws.Cells.Item(32, 5).Value = "  +1.29%  "
$ws.Cells.Item(33, 5).Value = "  -0.04%  "
$ws.Cells.Item(34, 5).Value = "  -0.20%  "
$ws.Cells.Item(35, 5).Value = "  +3.16%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "3.98"
$ws.Cells.Item(36, 5).Value = "  +0.77%  "
$ws.Cells.Item(37, 5).Value = "  -5.51%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "1.56"
$ws.Cells.Item(38, 5).Value = "  +3.01%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "38.58"
$ws.Cells.Item(39, 5).Value = "  +3.16%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "150.90"
$ws.Cells.Item(40, 5).Value = "  +9.57%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.375"
$ws.Cells.Item(41, 5).Value = "  -0.88%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "3.59"
$ws.Cells.Item(42, 5).Value = "  +1.31%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "281.63"
$ws.Cells.Item(43, 5).Value = "  +2.43%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "5.12"
$ws.Cells.Item(44, 5).Value = "  +1.79%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.0925"
$ws.Cells.Item(45, 5).Value = "  -0.24%  "
$ws.Cells.Item(46, 5).Value = "  -0.10%  "
$ws.Cells.Item(47, 5).Value = "  +0.03%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "18.33"
$ws.Cells.Item(48, 5).Value = "  +7.63%  "
$ws.Cells.Item(49, 5).Value = "  +0.21%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "17.09"
$ws.Cells.Item(50, 5).Value = "  +2.40%  "
$ws.Cells.Item(51, 5).Value = "  +1.26%  "
